$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Wagi_Model_Bez_Outlierow")
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(256, 70)'
$arr[0,1] = 17920
$arr[0,2] = -0.01711405254900455
$arr[0,3] = 0.1603449732065201
$arr[0,4] = -0.7438752055168152
$arr[0,5] = 0.695008397102356
$arr[0,6] = -0.008567610755562782
$arr[0,7] = -0.1068862080574036
$arr[0,8] = 0.07274356484413147
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 1.438883543014526
$arr[0,12] = '{''negative_ratio'': np.float64(0.5322544642857143), ''positive_ratio'': np.float64(0.46774553571428573), ''near_zero_ratio'': np.float64(0.5361049107142857)}'
$ws.Range("B2:N2").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(128, 256)'
$arr[0,1] = 32768
$arr[0,2] = -0.05280745029449463
$arr[0,3] = 0.1543912589550018
$arr[0,4] = -1.139596462249756
$arr[0,5] = 0.5488262176513672
$arr[0,6] = -0.01943532377481461
$arr[0,7] = -0.1248481050133705
$arr[0,8] = 0.04602491855621338
$arr[0,9] = 26
$arr[0,10] = 3
$arr[0,11] = 1.688422679901123
$arr[0,12] = '{''negative_ratio'': np.float64(0.57470703125), ''positive_ratio'': np.float64(0.42529296875), ''near_zero_ratio'': np.float64(0.597869873046875)}'
$ws.Range("B3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(64, 128)'
$arr[0,1] = 8192
$arr[0,2] = -0.08466237783432007
$arr[0,3] = 0.1818512678146362
$arr[0,4] = -1.142301201820374
$arr[0,5] = 0.7446489334106445
$arr[0,6] = -0.03681820631027222
$arr[0,7] = -0.2019456475973129
$arr[0,8] = 0.03562135994434357
$arr[0,9] = 1
$arr[0,10] = 1
$arr[0,11] = 1.886950135231018
$arr[0,12] = '{''negative_ratio'': np.float64(0.620361328125), ''positive_ratio'': np.float64(0.379638671875), ''near_zero_ratio'': np.float64(0.5172119140625)}'
$ws.Range("B4:N4").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(2, 64)'
$arr[0,1] = 128
$arr[0,2] = -0.0002197984140366316
$arr[0,3] = 0.1146299988031387
$arr[0,4] = -0.2582334876060486
$arr[0,5] = 0.2840132713317871
$arr[0,6] = 0.008058508858084679
$arr[0,7] = -0.1018093079328537
$arr[0,8] = 0.09006457775831223
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5422467589378357
$arr[0,12] = '{''negative_ratio'': np.float64(0.46875), ''positive_ratio'': np.float64(0.53125), ''near_zero_ratio'': np.float64(0.53125)}'
$ws.Range("B5:N5").Value = $arr

$ws = $wb.Worksheets.Item("Wagi_Model_Z_Outlierami")
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(256, 70)'
$arr[0,1] = 17920
$arr[0,2] = -0.008874831721186638
$arr[0,3] = 0.1233037412166595
$arr[0,4] = -0.7738926410675049
$arr[0,5] = 0.4860520660877228
$arr[0,6] = -0.004700297489762306
$arr[0,7] = -0.07732301950454712
$arr[0,8] = 0.05924554169178009
$arr[0,9] = 1
$arr[0,10] = 0
$arr[0,11] = 1.259944677352905
$arr[0,12] = '{''negative_ratio'': np.float64(0.52421875), ''positive_ratio'': np.float64(0.47578125), ''near_zero_ratio'': np.float64(0.6307477678571428)}'
$ws.Range("B2:N2").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(128, 256)'
$arr[0,1] = 32768
$arr[0,2] = -0.05243252962827682
$arr[0,3] = 0.1463489234447479
$arr[0,4] = -0.9223895072937012
$arr[0,5] = 0.5513448119163513
$arr[0,6] = -0.01403271965682507
$arr[0,7] = -0.1144021451473236
$arr[0,8] = 0.04133421182632446
$arr[0,9] = 79
$arr[0,10] = 0
$arr[0,11] = 1.473734378814697
$arr[0,12] = '{''negative_ratio'': np.float64(0.56787109375), ''positive_ratio'': np.float64(0.43212890625), ''near_zero_ratio'': np.float64(0.6468505859375)}'
$ws.Range("B3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(64, 128)'
$arr[0,1] = 8192
$arr[0,2] = -0.07360292971134186
$arr[0,3] = 0.16595458984375
$arr[0,4] = -0.7039116621017456
$arr[0,5] = 0.5289499759674072
$arr[0,6] = -0.03143126145005226
$arr[0,7] = -0.1800456941127777
$arr[0,8] = 0.03734163194894791
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 1.232861638069153
$arr[0,12] = '{''negative_ratio'': np.float64(0.60888671875), ''positive_ratio'': np.float64(0.39111328125), ''near_zero_ratio'': np.float64(0.5406494140625)}'
$ws.Range("B4:N4").Value = $arr
$arr = New-Object 'object[,]' 1,13
$arr[0,0] = '(2, 64)'
$arr[0,1] = 128
$arr[0,2] = -0.0008837929926812649
$arr[0,3] = 0.1321490854024887
$arr[0,4] = -0.3066065013408661
$arr[0,5] = 0.2701490223407745
$arr[0,6] = 0.01635450311005116
$arr[0,7] = -0.1083608865737915
$arr[0,8] = 0.1086637452244759
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5767555236816406
$arr[0,12] = '{''negative_ratio'': np.float64(0.46875), ''positive_ratio'': np.float64(0.53125), ''near_zero_ratio'': np.float64(0.4140625)}'
$ws.Range("B5:N5").Value = $arr

$ws = $wb.Worksheets.Item("Biasy_Model_Bez_Outlierow")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(256,)'
$arr[0,1] = 256
$arr[0,2] = -0.5612590909004211
$arr[0,3] = 0.2049462050199509
$arr[0,4] = -1.172370314598083
$arr[0,5] = -0.01639115251600742
$arr[0,6] = -0.5743221044540405
$ws.Range("B2:H2").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(128,)'
$arr[0,1] = 128
$arr[0,2] = 0.00926684308797121
$arr[0,3] = 0.3524489104747772
$arr[0,4] = -1.096748352050781
$arr[0,5] = 0.7321162819862366
$arr[0,6] = 0.03953400254249573
$ws.Range("B3:H3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(64,)'
$arr[0,1] = 64
$arr[0,2] = 0.445732593536377
$arr[0,3] = 0.08467815816402435
$arr[0,4] = 0.2514960467815399
$arr[0,5] = 0.6003814935684204
$arr[0,6] = 0.4430473446846008
$ws.Range("B4:H4").Value = $arr

$ws = $wb.Worksheets.Item("Biasy_Model_Z_Outlierami")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(256,)'
$arr[0,1] = 256
$arr[0,2] = -0.3838124871253967
$arr[0,3] = 0.134968176484108
$arr[0,4] = -0.7661592364311218
$arr[0,5] = 0.05474037304520607
$arr[0,6] = -0.3820775151252747
$ws.Range("B2:H2").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(128,)'
$arr[0,1] = 128
$arr[0,2] = 0.02567333169281483
$arr[0,3] = 0.2512891590595245
$arr[0,4] = -0.5473055839538574
$arr[0,5] = 0.7028363943099976
$arr[0,6] = 0.001539497869089246
$ws.Range("B3:H3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '(64,)'
$arr[0,1] = 64
$arr[0,2] = 0.3147733211517334
$arr[0,3] = 0.1046170219779015
$arr[0,4] = 0.02957485429942608
$arr[0,5] = 0.533165454864502
$arr[0,6] = 0.3192217648029327
$ws.Range("B4:H4").Value = $arr

$ws = $wb.Worksheets.Item("Biasy_Model_Bez_Outlierow")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = -0.08648733794689178
$arr[0,1] = 0.0126589797437191
$arr[0,2] = -0.09914632141590118
$arr[0,3] = -0.07382836192846298
$arr[0,4] = -0.08648733794689178
$ws.Range("D5:H5").Value = $arr

$ws = $wb.Worksheets.Item("Biasy_Model_Z_Outlierami")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = -0.1392880976200104
$arr[0,1] = 0.01767654344439507
$arr[0,2] = -0.1569646447896957
$arr[0,3] = -0.1216115579009056
$arr[0,4] = -0.1392880976200104
$ws.Range("D5:H5").Value = $arr

$ws = $wb.Worksheets.Item("Waznosc_Cech_Bez_Outlierow")
$arr = New-Object 'object[,]' 14,3
$arr[0,0] = 'data__tagData__linearAcceleration__y'
$arr[0,1] = 0.1594968289136887
$arr[0,2] = 0.7483857870101929
$arr[1,0] = 'data__tagData__quaternion__x'
$arr[1,1] = 0.1507901549339294
$arr[1,2] = 0.7075326442718506
$arr[2,0] = 'data__tagData__linearAcceleration__x'
$arr[2,1] = 0.1438106447458267
$arr[2,2] = 0.6747835874557495
$arr[3,0] = 'data__tagData__quaternion__y'
$arr[3,1] = 0.1433026343584061
$arr[3,2] = 0.6723999381065369
$arr[4,0] = 'data__tagData__magnetic__x'
$arr[4,1] = 0.1378122717142105
$arr[4,2] = 0.6466382145881653
$arr[5,0] = 'data__tagData__quaternion__z'
$arr[5,1] = 0.1320609897375107
$arr[5,2] = 0.619652271270752
$arr[6,0] = 'data__tagData__magnetic__z'
$arr[6,1] = 0.1298420280218124
$arr[6,2] = 0.6092405319213867
$arr[7,0] = 'data__tagData__quaternion__w'
$arr[7,1] = 0.1210319921374321
$arr[7,2] = 0.5679023265838623
$arr[8,0] = 'data__tagData__pressure'
$arr[8,1] = 0.1203977465629578
$arr[8,2] = 0.5649263858795166
$arr[9,0] = 'data__tagData__magnetic__y'
$arr[9,1] = 0.1117235496640205
$arr[9,2] = 0.5242255926132202
$arr[10,0] = 'data__tagData__gyro__z'
$arr[10,1] = 0.06795907020568848
$arr[10,2] = 0.3188753128051758
$arr[11,0] = 'data__tagData__linearAcceleration__z'
$arr[11,1] = 0.06748020648956299
$arr[11,2] = 0.3166284263134003
$arr[12,0] = 'data__tagData__gyro__y'
$arr[12,1] = 0.0470699667930603
$arr[12,2] = 0.2208601534366608
$arr[13,0] = 'data__tagData__gyro__x'
$arr[13,1] = 0.04119820147752762
$arr[13,2] = 0.1933088600635529
$ws.Range("A2:C15").Value = $arr

$ws = $wb.Worksheets.Item("Waznosc_Cech_Z_Outlierami")
$arr = New-Object 'object[,]' 14,3
$arr[0,0] = 'data__tagData__magnetic__z'
$arr[0,1] = 0.1268749982118607
$arr[0,2] = 0.7880849838256836
$arr[1,0] = 'data__tagData__quaternion__x'
$arr[1,1] = 0.1235200688242912
$arr[1,2] = 0.7672458291053772
$arr[2,0] = 'data__tagData__magnetic__x'
$arr[2,1] = 0.1187189370393753
$arr[2,2] = 0.7374235391616821
$arr[3,0] = 'data__tagData__quaternion__y'
$arr[3,1] = 0.1174513921141624
$arr[3,2] = 0.7295501828193665
$arr[4,0] = 'data__tagData__quaternion__w'
$arr[4,1] = 0.1046234965324402
$arr[4,2] = 0.6498696208000183
$arr[5,0] = 'data__tagData__pressure'
$arr[5,1] = 0.1025804728269577
$arr[5,2] = 0.6371793746948242
$arr[6,0] = 'data__tagData__linearAcceleration__y'
$arr[6,1] = 0.09987609833478928
$arr[6,2] = 0.6203811168670654
$arr[7,0] = 'data__tagData__magnetic__y'
$arr[7,1] = 0.09966656565666199
$arr[7,2] = 0.61907958984375
$arr[8,0] = 'data__tagData__quaternion__z'
$arr[8,1] = 0.09737355262041092
$arr[8,2] = 0.6048365235328674
$arr[9,0] = 'data__tagData__linearAcceleration__x'
$arr[9,1] = 0.08796892315149307
$arr[9,2] = 0.546419620513916
$arr[10,0] = 'data__tagData__linearAcceleration__z'
$arr[10,1] = 0.07794170081615448
$arr[10,2] = 0.4841354489326477
$arr[11,0] = 'data__tagData__gyro__z'
$arr[11,1] = 0.02600033208727837
$arr[11,2] = 0.1615012586116791
$arr[12,0] = 'data__tagData__gyro__y'
$arr[12,1] = 0.02569006010890007
$arr[12,2] = 0.1595740020275116
$arr[13,0] = 'data__tagData__gyro__x'
$arr[13,1] = 0.02475104667246342
$arr[13,2] = 0.1537413001060486
$ws.Range("A2:C15").Value = $arr

$ws = $wb.Worksheets.Item("Porownanie_Modeli")
$arr = New-Object 'object[,]' 1,9
$arr[0,0] = -0.01711405254900455
$arr[0,1] = -0.008874831721186638
$arr[0,2] = 0.008239220827817917
$arr[0,3] = 0.1603449732065201
$arr[0,4] = 0.1233037412166595
$arr[0,5] = -0.03704123198986053
$arr[0,6] = 1.438883543014526
$arr[0,7] = 1.259944677352905
$arr[0,8] = -0.1789388656616211
$ws.Range("B2:J2").Value = $arr
$arr = New-Object 'object[,]' 1,9
$arr[0,0] = -0.05280745029449463
$arr[0,1] = -0.05243252962827682
$arr[0,2] = 0.000374920666217804
$arr[0,3] = 0.1543912589550018
$arr[0,4] = 0.1463489234447479
$arr[0,5] = -0.008042335510253906
$arr[0,6] = 1.688422679901123
$arr[0,7] = 1.473734378814697
$arr[0,8] = -0.2146883010864258
$ws.Range("B3:J3").Value = $arr
$arr = New-Object 'object[,]' 1,9
$arr[0,0] = -0.08466237783432007
$arr[0,1] = -0.07360292971134186
$arr[0,2] = 0.01105944812297821
$arr[0,3] = 0.1818512678146362
$arr[0,4] = 0.16595458984375
$arr[0,5] = -0.01589667797088623
$arr[0,6] = 1.886950135231018
$arr[0,7] = 1.232861638069153
$arr[0,8] = -0.6540884971618652
$ws.Range("B4:J4").Value = $arr
$arr = New-Object 'object[,]' 1,9
$arr[0,0] = -0.0002197984140366316
$arr[0,1] = -0.0008837929926812649
$arr[0,2] = -0.0006639945786446333
$arr[0,3] = 0.1146299988031387
$arr[0,4] = 0.1321490854024887
$arr[0,5] = 0.01751908659934998
$arr[0,6] = 0.5422467589378357
$arr[0,7] = 0.5767555236816406
$arr[0,8] = 0.03450876474380493
$ws.Range("B5:J5").Value = $arr

Write-Host "done"